# Update from MV -datos- : append 20 new daily rows (03-08-2021 .. 30-08-2021)
# to the liquidity injection/drain table on Sheet1 (rows 151-170).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("03-08-2021", 38, 7, 10, -2, 0, -3, -8, -38, -11, 6),
    @("04-08-2021", 38, 7, 10, -2, 0, -2, -7, -37, -12, 6),
    @("05-08-2021", 38, 7, 10, -2, 0, -2, -6, -38, -11, 6),
    @("06-08-2021", 38, 7, 10, -2, 0, -2, -7, -38, -10, 6),
    @("09-08-2021", 37, 7, 10, -2, 0, -2, -2, -36, -17, 6),
    @("10-08-2021", 37, 7, 10, -2, 0, -2, -2, -35, -18, 6),
    @("11-08-2021", 37, 7, 10, -2, 0, -2, -3, -35, -18, 6),
    @("12-08-2021", 38, 7, 10, -2, 0, -2, -3, -35, -18, 6),
    @("13-08-2021", 38, 7, 10, -2, 0, -2, -3, -36, -17, 6),
    @("16-08-2021", 38, 7, 10, -2, 0, -2, -4, -35, -17, 6),
    @("17-08-2021", 37, 7, 10, -2, 0, -2, -4, -34, -17, 6),
    @("18-08-2021", 37, 7, 10, -2, 0, -2, -4, -34, -17, 6),
    @("19-08-2021", 37, 7, 10, -2, 0, -2, -5, -35, -15, 6),
    @("20-08-2021", 37, 7, 10, -2, 0, -2, -5, -35, -16, 6),
    @("23-08-2021", 37, 7, 10, -2, 0, -2, -5, -35, -15, 6),
    @("24-08-2021", 37, 7, 10, -2, 0, -2, -5, -36, -15, 6),
    @("25-08-2021", 37, 7, 10, -2, 0, -2, -4, -36, -16, 6),
    @("26-08-2021", 37, 7, 10, -2, 0, -2, -5, -35, -16, 6),
    @("27-08-2021", 37, 7, 10, -2, 0, -2, -6, -34, -15, 6),
    @("30-08-2021", 37, 7, 10, -2, 0, -2, -5, -34, -17, 6)
)

$startRow = 151
$endRow = $startRow + $data.Count - 1

# Column A holds dd-mm-yyyy text labels. Excel's COM layer auto-converts such
# strings into date serials unless the cell is already formatted as Text, so
# temporarily mark the target range as Text, enter the values, then restore
# the default "Normal" style (matching the rest of the column, which carries
# no explicit style) once the text values are locked in.
$colA = $ws.Range("A$startRow`:A$endRow")
$colA.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    for ($c = 1; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}

$colA.Style = "Normal"
